# Applies the "Added query to reset counter to 1 in ACCESS db" edit to the
# Tickets workbook. This updates the project ticket dates (shifted from the
# 3/27 milestone to 3/30), a couple of resource re-assignments, and the
# time/duration estimates for the final "9.x" submission tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Log")

# --- 1.5.1 Consolidate updates and changes: end date pushed out to 3/30 ---
$ws.Range("B12").Value = "2/10-03/30"

# --- 5.2 Create form: Back/Next button resource swapped from Gina to Mara ---
$ws.Range("D79").Value = "Mara"   # 5.2.2 Back button
$ws.Range("D81").Value = "Mara"   # 5.2.4 Next button

# --- 5.2.5 / 5.3.3 / 5.4.x: day range pushed from 3/21-3/27 to 3/21-3/30 ---
$ws.Range("B82").Value = "3/21-3/30"   # 5.2.5 Send input to database
$ws.Range("D82").Value = "Gina, Domi"  # resource gains Gina
$ws.Range("B86").Value = "3/21-3/30"   # 5.3.3 Confirm button
$ws.Range("B88").Value = "3/21-3/30"   # 5.4.1 Edit data
$ws.Range("B89").Value = "3/21-3/30"   # 5.4.2 Complete data
$ws.Range("D89").Value = "Gina, Domi"  # resource gains Gina
$ws.Range("B90").Value = "3/21-3/30"   # 5.4.3 Cancel data

# --- 6. Debugging: day range pushed from 3/20-3/27 to 3/20-3/30 ---
$ws.Range("B92").Value = "3/20-3/30"   # 6.1 Home form
$ws.Range("B93").Value = "3/20-3/30"   # 6.2 Create form
$ws.Range("B94").Value = "3/20-3/30"   # 6.3 Reservation form
$ws.Range("B95").Value = "3/20-3/30"   # 6.4 Records form

# --- 7. Quality Checking: day range pushed from 3/25-3/27 to 3/25-3/30 ---
$ws.Range("B97").Value = "3/25-3/30"    # 7.1 Check system flow
$ws.Range("B98").Value = "3/25-3/30"    # 7.2 Check syntax
$ws.Range("B99").Value = "3/25-3/30"    # 7.3 Check logical functionality
$ws.Range("B100").Value = "3/25-3/30"   # 7.4 Check system functionality
$ws.Range("B101").Value = "3/25-3/30"   # 7.5 Check user friendliness

# --- 8.3 / 9. dates pushed from 3/27/2018 (43186) to 3/30/2018 (43189) ---
$ws.Range("B105").Value = 43189   # 8.3 System Code
$ws.Range("B106").Value = 43189   # 9. Submission of Project
$ws.Range("B108").Value = 43189   # 9.2 Final commits using Git
$ws.Range("B109").Value = 43189   # 9.3 Push final commit on Github

# --- Updated duration estimates for the submission tasks ---
$ws.Range("C107").Value = "4 hrs"    # 9.1 Compilation of changes
$ws.Range("C108").Value = "15 mins"  # 9.2 Final commits using Git
$ws.Range("C109").Value = "15 mins"  # 9.3 Push final commit on Github

# --- Update the saved cursor/selection position ---
$ws.Activate()
$ws.Range("E89").Select()
